$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J, matching the style used by the other header cells (B1:H1)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Fill data rows 2-39: I = 1 (constant), J = copy of H's value
for ($r = 2; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
